$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold free-form text like "1.100"
# or "27.907.35" that Excel would otherwise auto-parse as numbers.
# Force text format while writing the new values, then restore the
# cells to the default (unstyled) look, matching the original
# formatting.
$range = $ws.Range("D2:E51")
$range.NumberFormat = "@"

$ws.Range("D2").Value = "27.907.35"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.812.54"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.4971"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("D8").Value = "0.3928"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("D9").Value = "0.09560"
$ws.Range("E9").Value = "  +22.93%  "
$ws.Range("D10").Value = "1.100"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "40.86"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "6.410"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "1.004"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "20.42"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "1.806.44"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "7.276"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "0.00001120"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").Value = "92.23"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "0.06658"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "17.12"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "5.909"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "27.957.04"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "2.255"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "159.27"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "2.023.49"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").Value = "20.53"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "2.384"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").Value = "127.63"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "1.033"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "5.560"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "3.621"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "0.06710"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("D36").Value = "8.927"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").Value = "0.02324"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "0.2134"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "4.936"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "11.19"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "0.6153"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "1.143"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "13.07"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "1.293"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").Value = "0.5875"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "3.699"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "122.88"
$ws.Range("E48").Value = "  -3.58%  "
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "1.176"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "0.06762"
$ws.Range("E51").Value = "  -0.49%  "

$range.Style = "Normal"
